# Tambah 1 column data test modul 12
# Adds a new "Nama Belakang" column (column C) with two data values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column C
$ws.Range("C1").Value = "Nama Belakang"

# New data values for column C
$ws.Range("C2").Value = "Hula"
$ws.Range("C4").Value = "HUla"

# Update the active selection to reflect where the user left off editing
[void]$ws.Range("C5").Select()
